$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 0.7170026666666667
$ws.Range("H2").Value2 = 2.151008
$ws.Range("I2").Value2 = 0.02953485643833859
$ws.Range("J2").Value2 = 0.02953485643833859
$ws.Range("M2").Value2 = 1.485259333333333
$ws.Range("N2").Value2 = 4.455778
$ws.Range("O2").Value2 = 0.3057455162066235
$ws.Range("P2").Value2 = 0.3057455162066235
$ws.Range("Q2").Value2 = 1.064934902691556
$ws.Range("R2").Value2 = 9.584414124224001
$ws.Range("S2").Value2 = 0.00903014992782835
$ws.Range("T2").Value2 = 0.00903014992782835
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 0.7170026666666667
$ws.Range("H3").Value2 = 2.151008
$ws.Range("I3").Value2 = 0.02953485643833859
$ws.Range("J3").Value2 = 0.02953485643833859
$ws.Range("O3").Value2 = 0.2805555239151429
$ws.Range("P3").Value2 = 0.2805555239151429
$ws.Range("Q3").Value2 = 0.9771962423751113
$ws.Range("R3").Value2 = 8.794766181376001
$ws.Range("S3").Value2 = 0.008286167121816615
$ws.Range("T3").Value2 = 0.008286167121816616
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.7170026666666667
$ws.Range("H4").Value2 = 2.151008
$ws.Range("I4").Value2 = 0.02953485643833859
$ws.Range("J4").Value2 = 0.02953485643833859
$ws.Range("O4").Value2 = 0.4136989598782336
$ws.Range("P4").Value2 = 0.4136989598782336
$ws.Range("Q4").Value2 = 1.440944963143111
$ws.Range("R4").Value2 = 12.968504668288
$ws.Range("S4").Value2 = 0.01221853938869363
$ws.Range("T4").Value2 = 0.01221853938869363
$ws.Range("I5").Value2 = 0.4970672037825566
$ws.Range("J5").Value2 = 0.4970672037825566
$ws.Range("M5").Value2 = 1.485259333333333
$ws.Range("N5").Value2 = 4.455778
$ws.Range("O5").Value2 = 0.3057455162066235
$ws.Range("P5").Value2 = 0.3057455162066235
$ws.Range("Q5").Value2 = 17.92269467760845
$ws.Range("R5").Value2 = 161.304252098476
$ws.Range("S5").Value2 = 0.1519760688098807
$ws.Range("T5").Value2 = 0.1519760688098807
$ws.Range("I6").Value2 = 0.4970672037825566
$ws.Range("J6").Value2 = 0.4970672037825566
$ws.Range("O6").Value2 = 0.2805555239151429
$ws.Range("P6").Value2 = 0.2805555239151429
$ws.Range("S6").Value2 = 0.1394549497782503
$ws.Range("T6").Value2 = 0.1394549497782503
$ws.Range("I7").Value2 = 0.4970672037825566
$ws.Range("J7").Value2 = 0.4970672037825566
$ws.Range("O7").Value2 = 0.4136989598782336
$ws.Range("P7").Value2 = 0.4136989598782336
$ws.Range("Q7").Value2 = 24.25088759545688
$ws.Range("S7").Value2 = 0.2056361851944257
$ws.Range("T7").Value2 = 0.2056361851944257
$ws.Range("I8").Value2 = 0.4733979397791048
$ws.Range("J8").Value2 = 0.4733979397791048
$ws.Range("M8").Value2 = 1.485259333333333
$ws.Range("N8").Value2 = 4.455778
$ws.Range("O8").Value2 = 0.3057455162066235
$ws.Range("P8").Value2 = 0.3057455162066235
$ws.Range("Q8").Value2 = 17.06925476294623
$ws.Range("R8").Value2 = 153.623292866516
$ws.Range("S8").Value2 = 0.1447392974689145
$ws.Range("T8").Value2 = 0.1447392974689145
$ws.Range("I9").Value2 = 0.4733979397791048
$ws.Range("J9").Value2 = 0.4733979397791048
$ws.Range("O9").Value2 = 0.2805555239151429
$ws.Range("P9").Value2 = 0.2805555239151429
$ws.Range("S9").Value2 = 0.132814407015076
$ws.Range("T9").Value2 = 0.132814407015076
$ws.Range("I10").Value2 = 0.4733979397791048
$ws.Range("J10").Value2 = 0.4733979397791048
$ws.Range("O10").Value2 = 0.4136989598782336
$ws.Range("P10").Value2 = 0.4136989598782336
$ws.Range("S10").Value2 = 0.1958442352951144
$ws.Range("T10").Value2 = 0.1958442352951144
